# 11.5.1 — add a new "2023" data column (T) to the right of the existing
# "2022" column (S), carrying over each row's number format, and bump a
# handful of 2019 (column P) figures that were revised upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column geometry -------------------------------------------------
# Columns A:C become a uniform width, D:T (now including the new T) get
# the narrower "year" column width.
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 31.75
$ws.Range("D1:T1").EntireColumn.ColumnWidth = 7.59

# --- Row 1 header gets a touch more height for the wrapped title -----
$ws.Rows("1").RowHeight = 31.5

# --- Revised 2019 (column P) figures ----------------------------------
$ws.Range("P5").Value = 19
$ws.Range("P6").Value = 10
$ws.Range("P14").Value = 7
$ws.Range("P15").Value = 4

# --- New column T: clone each row's column-S formatting, then fill in
#     the 2023 values (header row + 30 data rows) -----------------------
$ws.Range("S3").Copy()
$ws.Range("T3").PasteSpecial(-4122)

$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = 2023

$values = @{
    5  = 44
    6  = 24
    7  = 20
    8  = "-"
    9  = "-"
    10 = "-"
    11 = 5
    12 = 1
    13 = 4
    14 = 8
    15 = 6
    16 = 2
    17 = 5
    18 = 1
    19 = 4
    20 = 7
    21 = 5
    22 = 2
    23 = "-"
    24 = "-"
    25 = "-"
    26 = 18
    27 = 10
    28 = 8
    29 = "-"
    30 = "-"
    31 = "-"
    32 = 1
    33 = 1
    34 = "-"
}

foreach ($r in 5..34) {
    $ws.Range("S$r").Copy()
    $ws.Range("T$r").PasteSpecial(-4122)
    $ws.Range("T$r").Value = $values[$r]
}

# --- Selection reverts to the default top-left cell -------------------
$ws.Range("A1").Select()
